$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8. This shifts the existing rows 8-53
# down to 9-54, carrying along their values and formatting (including
# the date style on column D).
$ws.Rows("8:8").Insert()

# Populate the newly inserted row 8 with the new data record.
$ws.Range("A8").Value = 5
$ws.Range("B8").Value = "Macroferia Regional de Talca"
$ws.Range("C8").Value = "Maule"
$ws.Range("D8").Value = 44503
$ws.Range("E8").Value = 7
$ws.Range("F8").Value = 100112022
$ws.Range("G8").Value = "Arveja Verde"
$ws.Range("H8").Value = "Sin especificar"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 600
$ws.Range("K8").Value = 10000
$ws.Range("L8").Value = 12000
$ws.Range("M8").Value = 11000
$ws.Range("N8").Value = "$/saco 25 kilos"
$ws.Range("O8").Value = "Región del Maule"
$ws.Range("P8").Value = 440
$ws.Range("Q8").Value = 25
$ws.Range("R8").Value = "Hortaliza"
